# Slide 1, "Subtitle 2" placeholder: second paragraph holds the phone
# number "2034757030". The author inserted an extra "2" right after the
# leading "20", turning it into "20234757030".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$numberParagraph = $tr.Paragraphs(2, 1)
$numberRun = $numberParagraph.Runs(1, 1)
$numberRun.Text = "20234757030"
